$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, matching the
# original inline-string cell type (values like "1.00" or "61.701.19"
# must not be reinterpreted as numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.701.19'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.450.19'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.69'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.91'
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.585'
$ws.Range("E8").Value = '  -0.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.449.22'
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("E10").Value = '  +1.88%  '

$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.43'
$ws.Range("E12").Value = '  +3.23%  '

$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.98'
$ws.Range("E14").Value = '  +2.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.894.34'
$ws.Range("E15").Value = '  -0.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000169'
$ws.Range("E16").Value = '  +5.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.784.21'
$ws.Range("E17").Value = '  +0.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.452.87'
$ws.Range("E18").Value = '  -0.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.81'
$ws.Range("E19").Value = '  -1.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.95'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '319.54'
$ws.Range("E22").Value = '  +1.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.86'
$ws.Range("E24").Value = '  +8.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.83'
$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0975'
$ws.Range("E26").Value = '  -2.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.572.07'

$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.81'
$ws.Range("E30").Value = '  +5.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '530.99'
$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -0.20%  '

$ws.Range("E33").Value = '  +2.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  +3.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.66'
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("E39").Value = '  +2.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.25'
$ws.Range("E40").Value = '  +0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  +5.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '139.35'
$ws.Range("E42").Value = '  -1.15%  '

$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.27'
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("E45").Value = '  +1.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.76'
$ws.Range("E46").Value = '  -2.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.61'
$ws.Range("E47").Value = '  +3.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.30'
$ws.Range("E48").Value = '  +2.95%  '

$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("E50").Value = '  +1.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0933'
$ws.Range("E51").Value = '  +0.13%  '
